# Updated cryptos list on Wed Jan 17 02:20:28 UTC 2024 with GitHub Actions
# Applies the latest coinranking.com snapshot values to the sheet while
# preserving each cell's original (unformatted, text) storage type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $cell = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "315.46") are not
    # silently reinterpreted by Excel as numbers; then drop back to the
    # workbook's default "Normal" style so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

Set-TextValue "D2" '42.999.60'
Set-TextValue "E2" '  +0.82%  '
Set-TextValue "D3" '2.578.72'
Set-TextValue "E3" '  +2.37%  '
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "D5" '315.46'
Set-TextValue "E5" '  -0.80%  '
Set-TextValue "D6" '99.17'
Set-TextValue "E6" '  +3.72%  '
Set-TextValue "D7" '0.576'
Set-TextValue "E7" '  +0.20%  '
Set-TextValue "D9" '0.539'
Set-TextValue "E9" '  +1.29%  '
Set-TextValue "D10" '36.08'
Set-TextValue "E10" '  +0.87%  '
Set-TextValue "D11" '0.0813'
Set-TextValue "E11" '  +0.25%  '
Set-TextValue "D12" '7.58'
Set-TextValue "E12" '  +0.32%  '
Set-TextValue "D13" '2.977.34'
Set-TextValue "E13" '  +2.48%  '
Set-TextValue "E14" '  -0.12%  '
Set-TextValue "D15" '2.694.42'
Set-TextValue "E15" '  +6.57%  '
Set-TextValue "D16" '15.69'
Set-TextValue "E16" '  +2.48%  '
Set-TextValue "D17" '0.843'
Set-TextValue "E17" '  -1.00%  '
Set-TextValue "D18" '43.077.43'
Set-TextValue "E18" '  +0.86%  '
Set-TextValue "E19" '  +2.51%  '
Set-TextValue "D20" '12.73'
Set-TextValue "E20" '  -0.97%  '
Set-TextValue "D21" '0.0₃0971'
Set-TextValue "E21" '  +1.25%  '
Set-TextValue "D22" '69.53'
Set-TextValue "E22" '  -0.19%  '
Set-TextValue "D23" '250.56'
Set-TextValue "E23" '  +0.32%  '
Set-TextValue "D24" '2.96'
Set-TextValue "E24" '  +0.22%  '
Set-TextValue "D25" '2.10'
Set-TextValue "E25" '  -0.25%  '
Set-TextValue "D26" '27.13'
Set-TextValue "E26" '  +2.10%  '
Set-TextValue "E27" '  -0.03%  '
Set-TextValue "D28" '2.40'
Set-TextValue "E28" '  -0.90%  '
Set-TextValue "D29" '40.41'
Set-TextValue "E29" '  -2.79%  '
Set-TextValue "E30" '  +0.40%  '
Set-TextValue "D31" '5.86'
Set-TextValue "E31" '  -1.73%  '
Set-TextValue "D32" '158.07'
Set-TextValue "E32" '  +0.10%  '
Set-TextValue "D33" '3.45'
Set-TextValue "E33" '  +5.34%  '
Set-TextValue "D34" '2.14'
Set-TextValue "E34" '  +0.04%  '
Set-TextValue "D35" '0.0806'
Set-TextValue "E35" '  +3.53%  '
Set-TextValue "B36" 'WEMIXToken'
Set-TextValue "C36" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D36" '2.68'
Set-TextValue "E36" '  +0.19%  '
Set-TextValue "B37" 'Celestia'
Set-TextValue "C37" 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D37" '18.84'
Set-TextValue "E37" '  -2.26%  '
Set-TextValue "E38" '  +1.30%  '
Set-TextValue "D40" '24.48'
Set-TextValue "E40" '  +4.39%  '
Set-TextValue "E41" '  +0.62%  '
Set-TextValue "E42" '  +6.57%  '
Set-TextValue "E43" '  +0.19%  '
Set-TextValue "E44" '  -0.05%  '
Set-TextValue "D45" '3.26'
Set-TextValue "E45" '  -1.74%  '
Set-TextValue "D46" '2.009.99'
Set-TextValue "E46" '  -1.01%  '
Set-TextValue "D47" '8.93'
Set-TextValue "E47" '  +0.06%  '
Set-TextValue "B48" 'RocketPoolETH'
Set-TextValue "C48" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D48" '2.828.45'
Set-TextValue "E48" '  +2.64%  '
Set-TextValue "B49" 'Algorand'
Set-TextValue "C49" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D49" '0.197'
Set-TextValue "E49" '  +1.78%  '
Set-TextValue "B50" 'BitcoinSV'
Set-TextValue "C50" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D50" '82.39'
Set-TextValue "E50" '  -2.52%  '
Set-TextValue "B51" 'ordi'
Set-TextValue "C51" 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue "D51" '75.07'
Set-TextValue "E51" '  -0.34%  '
